$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 362.12766
$ws.Range("J17").Value = 362.12766
$ws.Range("L17").Value = 1086.38298
$ws.Range("N17").Value = -1422.38298
$ws.Range("H19").Value = 3971.9
$ws.Range("I19").Value = 12822.5
$ws.Range("J19").Value = 753.5
$ws.Range("K19").Value = 12822.5
$ws.Range("L19").Value = 753.5
$ws.Range("M19").Value = -12647.5
$ws.Range("N19").Value = -1103.5
$ws.Range("H40").Value = 1161.9324
$ws.Range("J40").Value = 1183.3273
$ws.Range("L40").Value = 1183.3273
$ws.Range("N40").Value = -1533.3273
$ws.Range("H64").Value = 4875.696
$ws.Range("I64").Value = 4452.231
$ws.Range("J64").Value = 5426.2
$ws.Range("K64").Value = 4452.231
$ws.Range("L64").Value = 5426.2
$ws.Range("M64").Value = -4204.231
$ws.Range("N64").Value = -5922.2
$ws.Range("H67").Value = 4875.696
$ws.Range("I67").Value = 4452.231
$ws.Range("J67").Value = 5426.2
$ws.Range("K67").Value = 4452.231
$ws.Range("L67").Value = 5426.2
$ws.Range("M67").Value = -3594.231
$ws.Range("N67").Value = -7142.2
$ws.Range("H116").Value = 3074.12
$ws.Range("I116").Value = 3136.6924
$ws.Range("J116").Value = 3006.3333
$ws.Range("K116").Value = 3136.6924
$ws.Range("L116").Value = 3006.3333
$ws.Range("M116").Value = 305.3076000000001
$ws.Range("N116").Value = -9890.3333
$ws.Range("H132").Value = 3127.4075
$ws.Range("I132").Value = 1506.3636
$ws.Range("K132").Value = 4519.0908
$ws.Range("M132").Value = -1989.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1076.9
$ws.Range("I45").Value = 923.0909
$ws.Range("J45").Value = 1499.875
$ws.Range("K45").Value = 923.0909
$ws.Range("L45").Value = 1499.875
$ws.Range("M45").Value = -546.0909
$ws.Range("N45").Value = -2253.875
$ws.Range("H61").Value = 3023.1592
$ws.Range("I61").Value = 2157.7
$ws.Range("K61").Value = 2157.7
$ws.Range("M61").Value = -1945.7
$ws.Range("H63").Value = 5166.6665
$ws.Range("I63").Value = 6500
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 6500
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -5814
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 5166.6665
$ws.Range("I66").Value = 6500
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 32500
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -29068
$ws.Range("N66").Value = -19364
$ws.Range("H132").Value = 27434.293
$ws.Range("I132").Value = 38764.965
$ws.Range("J132").Value = 3029.7693
$ws.Range("K132").Value = 116294.895
$ws.Range("L132").Value = 9089.3079
$ws.Range("M132").Value = -113764.895
$ws.Range("N132").Value = -14149.3079
$ws.Range("H136").Value = 3023.1592
$ws.Range("I136").Value = 2157.7
$ws.Range("K136").Value = 6473.099999999999
$ws.Range("M136").Value = -3923.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5541404
$ws.Range("I99").Value = 2027797.1
$ws.Range("K99").Value = 2027797.1
$ws.Range("M99").Value = -2026299.1
$ws.Range("H107").Value = 1796.75
$ws.Range("I107").Value = 2547.75
$ws.Range("J107").Value = 1045.75
$ws.Range("K107").Value = 2547.75
$ws.Range("L107").Value = 1045.75
$ws.Range("M107").Value = -627.75
$ws.Range("N107").Value = -4885.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2144.4038
$ws.Range("I31").Value = 1248.0857
$ws.Range("J31").Value = 3989.7646
$ws.Range("K31").Value = 1248.0857
$ws.Range("L31").Value = 3989.7646
$ws.Range("M31").Value = -953.0857000000001
$ws.Range("N31").Value = -4579.7646
$ws.Range("H34").Value = 2144.4038
$ws.Range("I34").Value = 1248.0857
$ws.Range("J34").Value = 3989.7646
$ws.Range("K34").Value = 1248.0857
$ws.Range("L34").Value = 3989.7646
$ws.Range("M34").Value = -1046.0857
$ws.Range("N34").Value = -4393.7646
$ws.Range("H58").Value = 1832.0322
$ws.Range("I58").Value = 932.86664
$ws.Range("K58").Value = 932.86664
$ws.Range("M58").Value = -729.86664
$ws.Range("H99").Value = 52112.75
$ws.Range("I99").Value = 112638.664
$ws.Range("K99").Value = 112638.664
$ws.Range("M99").Value = -111140.664
$ws.Range("H105").Value = 978.23914
$ws.Range("I105").Value = 859.3333
$ws.Range("J105").Value = 1406.3
$ws.Range("K105").Value = 859.3333
$ws.Range("L105").Value = 1406.3
$ws.Range("M105").Value = 887.6667
$ws.Range("N105").Value = -4900.3
$ws.Range("H106").Value = 24900
$ws.Range("J106").Value = 24900
$ws.Range("L106").Value = 24900
$ws.Range("N106").Value = -27424
$ws.Range("H126").Value = 52112.75
$ws.Range("I126").Value = 112638.664
$ws.Range("K126").Value = 337915.992
$ws.Range("M126").Value = -335445.992
$ws.Range("H132").Value = 1672.0256
$ws.Range("I132").Value = 1041.9032
$ws.Range("K132").Value = 3125.7096
$ws.Range("M132").Value = -595.7096000000001
$ws.Range("H134").Value = 1406.75
$ws.Range("I134").Value = 1009.125
$ws.Range("K134").Value = 3027.375
$ws.Range("M134").Value = -492.375
$ws.Range("H136").Value = 1832.0322
$ws.Range("I136").Value = 932.86664
$ws.Range("K136").Value = 2798.59992
$ws.Range("M136").Value = -248.5999199999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 527.7143
$ws.Range("I107").Value = 105.111115
$ws.Range("J107").Value = 844.6667
$ws.Range("K107").Value = 315.333345
$ws.Range("L107").Value = 2534.0001
$ws.Range("M107").Value = 1604.666655
$ws.Range("N107").Value = -6374.0001
$ws.Range("H131").Value = 1803.942
$ws.Range("I131").Value = 1300
$ws.Range("J131").Value = 1899.5172
$ws.Range("K131").Value = 3900
$ws.Range("L131").Value = 5698.5516
$ws.Range("M131").Value = 1140
$ws.Range("N131").Value = -15778.5516
$ws.Range("H132").Value = 5345.1577
$ws.Range("I132").Value = 2278
$ws.Range("J132").Value = 9562.5
$ws.Range("K132").Value = 20502
$ws.Range("L132").Value = 86062.5
$ws.Range("M132").Value = -17972
$ws.Range("N132").Value = -91122.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4782.353
$ws.Range("I70").Value = 4180.769
$ws.Range("J70").Value = 6737.5
$ws.Range("K70").Value = 4180.769
$ws.Range("L70").Value = 6737.5
$ws.Range("M70").Value = -3910.769
$ws.Range("N70").Value = -7277.5
$ws.Range("H73").Value = 4782.353
$ws.Range("I73").Value = 4180.769
$ws.Range("J73").Value = 6737.5
$ws.Range("K73").Value = 4180.769
$ws.Range("L73").Value = 6737.5
$ws.Range("M73").Value = -3244.769
$ws.Range("N73").Value = -8609.5
$ws.Range("H113").Value = 2065.2856
$ws.Range("I113").Value = 1802.75
$ws.Range("J113").Value = 2415.3333
$ws.Range("K113").Value = 1802.75
$ws.Range("L113").Value = 2415.3333
$ws.Range("M113").Value = 367.25
$ws.Range("N113").Value = -6755.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6426.263
$ws.Range("I40").Value = 7047
$ws.Range("J40").Value = 1150
$ws.Range("K40").Value = 7047
$ws.Range("L40").Value = 1150
$ws.Range("M40").Value = -6911
$ws.Range("N40").Value = -1422
